# Update the "repaymentstrategy" row on the ProductLoanInput sheet with the
# new scenario value and highlight it with the same green fill used
# elsewhere on the sheet, then leave that cell selected/active.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$cell = $ws.Range("B17")
$cell.Value = "Penalties, Fees, Interest, Principal order"
$cell.Interior.Color = 5296274
$cell.HorizontalAlignment = -4131
$cell.VerticalAlignment = -4160

$ws.Activate()
$cell.Select()
